# Update "想去人数" (want-to-go count) figures in the F column across sheets,
# reflecting refreshed scrape numbers baked into the gh-pages data output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        5  = 180
        6  = 1065
        7  = 1041
        8  = 8107
        10 = 201
        11 = 6864
        14 = 4943
        16 = 5369
        17 = 1072
        18 = 326
        19 = 330
        20 = 453
        26 = 9091
        28 = 1640
        29 = 585
        32 = 847
        37 = 1860
        40 = 4740
        43 = 70
        44 = 145
        46 = 34
        47 = 912
        48 = 1246
    }
    "演出" = @{
        6 = 23
        9 = 180
    }
    "全部类型" = @{
        6  = 180
        8  = 1065
        9  = 1041
        10 = 8107
        12 = 201
        13 = 6864
        17 = 4943
        19 = 5369
        20 = 1072
        21 = 326
        22 = 330
        23 = 453
        26 = 180
        27 = 9091
        29 = 1640
        30 = 585
        33 = 847
        38 = 1860
        41 = 4740
        44 = 70
        45 = 145
        47 = 912
        48 = 1246
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
